$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.145.22'
$ws.Range('E2').Value = '  -0.47%  '
$ws.Range('D3').Value = '2.269.46'
$ws.Range('E3').Value = '  +0.17%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '497.77'
$ws.Range('E5').Value = '  +0.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '128.24'
$ws.Range('E6').Value = '  +0.77%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  -0.63%  '
$ws.Range('E9').Value = '  +0.14%  '
$ws.Range('E10').Value = '  +0.42%  '
$ws.Range('E11').Value = '  +3.22%  '
$ws.Range('E12').Value = '  +1.91%  '
$ws.Range('D13').Value = '2.671.77'
$ws.Range('E13').Value = '  -1.13%  '
$ws.Range('E14').Value = '  +3.79%  '
$ws.Range('D15').Value = '54.115.49'
$ws.Range('E15').Value = '  -0.79%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000130'
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('D17').Value = '2.282.36'
$ws.Range('E17').Value = '  -0.66%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.21'
$ws.Range('E18').Value = '  +1.80%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.13'
$ws.Range('E19').Value = '  +1.77%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '302.96'
$ws.Range('E20').Value = '  -0.60%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.32'
$ws.Range('E21').Value = '  -2.04%  '
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '61.06'
$ws.Range('E23').Value = '  -3.14%  '
$ws.Range('E24').Value = '  -0.32%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.150'
$ws.Range('E25').Value = '  -1.44%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.29'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '170.57'
$ws.Range('E27').Value = '  -0.30%  '
$ws.Range('E28').Value = '  +0.26%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0₃0688'
$ws.Range('E29').Value = '  +0.39%  '
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.90'
$ws.Range('E30').Value = '  +0.22%  '
$ws.Range('E31').Value = '  +0.41%  '
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '17.72'
$ws.Range('E33').Value = '  +0.67%  '
$ws.Range('E34').Value = '  +0.46%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.946'
$ws.Range('E35').Value = '  +9.46%  '
$ws.Range('E36').Value = '  -1.30%  '
$ws.Range('E37').Value = '  +0.87%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.371'
$ws.Range('E38').Value = '  -0.94%  '
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('E40').Value = '  +0.39%  '
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '124.86'
$ws.Range('E42').Value = '  -3.29%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0491'
$ws.Range('E43').Value = '  +1.55%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0889'
$ws.Range('E44').Value = '  -0.62%  '
$ws.Range('E45').Value = '  -1.12%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '238.75'
$ws.Range('E46').Value = '  -1.52%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.370'
$ws.Range('E47').Value = '  -1.00%  '
$ws.Range('E48').Value = '  +0.65%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.75'
$ws.Range('E49').Value = '  +0.44%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '16.12'
$ws.Range('E50').Value = '  -1.50%  '
$ws.Range('E51').Value = '  -0.64%  '
